# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets,
# matching the regenerated data published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# row => [old, new] for sheet "展览" (sheet1)
$updatesExhibit = @{
    2  = 1337
    3  = 1221
    4  = 14645
    5  = 17859
    18 = 45
    19 = 1342
    20 = 149
    23 = 216
    24 = 7372
    27 = 42
    28 = 1183
    30 = 5878
    31 = 72
    32 = 49
    33 = 145
    36 = 5127
}

# row => new value for sheet "全部类型" (sheet4), row numbers shifted by +1
# starting from row 22 because that sheet has one extra data row (r22)
$updatesAll = @{
    2  = 1337
    3  = 1221
    4  = 14645
    5  = 17859
    18 = 45
    19 = 1342
    20 = 149
    24 = 216
    25 = 7372
    28 = 42
    29 = 1183
    32 = 5878
    33 = 72
    34 = 49
    35 = 145
    38 = 5127
}

$wsExhibit = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibit.Keys) {
    $wsExhibit.Range("F$row").Value = $updatesExhibit[$row]
}

$wsAll = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAll.Keys) {
    $wsAll.Range("F$row").Value = $updatesAll[$row]
}
